$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove existing hyperlinks (engine only supports bulk clear reliably) ---
$ws.Hyperlinks.Delete()

# --- Drop the old rows 8-17; only 6 data rows remain (rows 2-7) ---
$ws.Range("A8:H17").EntireRow.Delete()

# --- Row 2 ---
$ws.Cells(2,1).Value = "2026-02-14 06:41:21"
$ws.Cells(2,2).Value = "ruby使用のX,instagram,threadsなどのスクレイピングについて"
$ws.Cells(2,3).Value = "システム開発"
$ws.Cells(2,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells(2,5).Value = "期限情報なし"
$ws.Cells(2,6).Value = "https://www.lancers.jp/work/detail/5491704"
$ws.Cells(2,7).Value = 40
$ws.Cells(2,8).Value = "◆スクレイピング"

# --- Row 3 ---
$ws.Cells(3,1).Value = "2026-02-14 06:41:21"
$ws.Cells(3,2).Value = "スクレイピング高速化サポート募集"
$ws.Cells(3,3).Value = "システム開発"
$ws.Cells(3,4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells(3,5).Value = "期限情報なし"
$ws.Cells(3,6).Value = "https://www.lancers.jp/work/detail/5491672"
$ws.Cells(3,7).Value = 40
$ws.Cells(3,8).Value = "◆スクレイピング"

# --- Row 4 ---
$ws.Cells(4,1).Value = "2026-02-14 06:41:21"
$ws.Cells(4,2).Value = "bubbleで構築したサイトの修正対応"
$ws.Cells(4,3).Value = "システム開発"
$ws.Cells(4,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells(4,5).Value = "期限情報なし"
$ws.Cells(4,6).Value = "https://www.lancers.jp/work/detail/5491578"
$ws.Cells(4,7).Value = 30
$ws.Cells(4,8).Value = "◇サイト"

# --- Row 5 ---
$ws.Cells(5,1).Value = "2026-02-14 06:41:21"
$ws.Cells(5,2).Value = "bubbleで構築したサイトの修正対応"
$ws.Cells(5,3).Value = "システム開発"
$ws.Cells(5,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells(5,5).Value = "期限情報なし"
$ws.Cells(5,6).Value = "https://www.lancers.jp/work/detail/5491569"
$ws.Cells(5,7).Value = 30
$ws.Cells(5,8).Value = "◇サイト"

# --- Row 6 ---
$ws.Cells(6,1).Value = "2026-02-14 06:41:21"
$ws.Cells(6,2).Value = "【急募】Githubとロリポップサーバーの連携不具合解決依頼"
$ws.Cells(6,3).Value = "システム開発"
$ws.Cells(6,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells(6,5).Value = "期限情報なし"
$ws.Cells(6,6).Value = "https://www.lancers.jp/work/detail/5491736"
$ws.Cells(6,7).Value = 10
$ws.Cells(6,8).ClearContents()

# --- Row 7 ---
$ws.Cells(7,1).Value = "2026-02-14 06:41:21"
$ws.Cells(7,2).Value = "【3,000円 / 急募】GitHubとVercelの連携設定エラーの解消"
$ws.Cells(7,3).Value = "システム開発"
$ws.Cells(7,4).Value = "1,000 ~ 5,000 円 / 固定"
$ws.Cells(7,5).Value = "期限情報なし"
$ws.Cells(7,6).Value = "https://www.lancers.jp/work/detail/5491643"
$ws.Cells(7,7).Value = 10
$ws.Cells(7,8).ClearContents()

# --- Recreate hyperlinks for F2:F7 and restore the Hyperlink cell style ---
$urls = @(
    "https://www.lancers.jp/work/detail/5491704",
    "https://www.lancers.jp/work/detail/5491672",
    "https://www.lancers.jp/work/detail/5491578",
    "https://www.lancers.jp/work/detail/5491569",
    "https://www.lancers.jp/work/detail/5491736",
    "https://www.lancers.jp/work/detail/5491643"
)
for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws.Hyperlinks.Add($ws.Cells($r,6), $urls[$i])
    $ws.Cells($r,6).Style = "Hyperlink"
}

# --- Column widths (ColumnWidth units differ from stored "width" by 5/6) ---
$ws.Columns.Item(2).ColumnWidth = 42 - 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 26 - 0.8333333333333334
$ws.Columns.Item(8).ColumnWidth = 12 - 0.8333333333333334
